$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 24-45: rnaDate (col A) and s1cDNADate (col D) get their TDY date
# bumped from "01.09.17" to "01.09.18". These are plain text values (not
# real Excel dates) in the source sheet, so force text formatting before
# the write to stop Excel from auto-converting the string into a date
# serial, then restore the default "Normal" style so no stray per-cell
# number format is left behind.
$rangeA = $ws.Range("A24:A45")
$rangeD = $ws.Range("D24:D45")

$rangeA.NumberFormat = "@"
$rangeD.NumberFormat = "@"

for ($r = 24; $r -le 45; $r++) {
    $ws.Range("A$r").Value = "01.09.18"
    $ws.Range("D$r").Value = "01.09.18"
}

$rangeA.NumberFormat = "General"
$rangeD.NumberFormat = "General"
$rangeA.Style = "Normal"
$rangeD.Style = "Normal"

# Mirror the author's final selection (scrolled down toward the rows that
# were just edited).
$ws.Range("D24:D45").Select()
